{"js": "// Apply the \"Added many more features\" edit:\n//  - Update the title (appears twice: H1 heading + bold run near the end)\n//  - Update the meta description (italic run at the end)\n//  - Split/replace the \"What we like\" / \"What we don't like\" bullet items\nconst body = context.document.body;\n\nconst replacements = [\n  // Title appears twice verbatim (Heading1 + bold run) - both get the same text.\n  [\n    \"Play Krakatoa Lucky Break for Free - Game Review\",\n    \"Play Krakatoa Lucky Break Free Slot: Volatile Wins & Big Jackpots!\",\n  ],\n  // \"What we like\" bullets\n  [\n    \"Fair level of volatility with frequent smaller wins and rare big wins\",\n    \"Fair level of volatility\",\n  ],\n  [\n    \"Bonus feature with multipliers and jackpots up to 15,000 game credits\",\n    \"Frequent smaller wins\",\n  ],\n  [\n    \"Electrifying graphics and a vintage look with vivid colors\",\n    \"Rare big wins with jackpots\",\n  ],\n  [\n    \"Symbols related to volcanic eruptions add to the theme's allure\",\n    \"Electrifying graphics and vintage look\",\n  ],\n  // \"What we don't like\" bullets\n  [\n    \"High level of volatility may not appeal to everyone\",\n    \"High level of volatility\",\n  ],\n  [\n    \"30 paylines may seem limited compared to some other online slots\",\n    \"Could offer more paylines\",\n  ],\n  // Meta description (italic run)\n  [\n    \"Discover Krakatoa Lucky Break, an online slot with frequent wins and huge jackpots. Play for free and read our review of the game's mechanics.\",\n    \"Read our review of Krakatoa Lucky Break, a slot game with high volatility and big jackpots. Play for free!\",\n  ],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit:\n#  - Update the title (appears twice: H1 heading + bold run near the end)\n#  - Update the meta description (italic run at the end)\n#  - Split/replace the \"What we like\" / \"What we don't like\" bullet items\n$d = $word.ActiveDocument\n\n$replacements = @(\n    # Title appears twice verbatim (Heading1 + bold run) - both get the same text.\n    @(\"Play Krakatoa Lucky Break for Free - Game Review\", \"Play Krakatoa Lucky Break Free Slot: Volatile Wins & Big Jackpots!\"),\n    # \"What we like\" bullets\n    @(\"Fair level of volatility with frequent smaller wins and rare big wins\", \"Fair level of volatility\"),\n    @(\"Bonus feature with multipliers and jackpots up to 15,000 game credits\", \"Frequent smaller wins\"),\n    @(\"Electrifying graphics and a vintage look with vivid colors\", \"Rare big wins with jackpots\"),\n    @(\"Symbols related to volcanic eruptions add to the theme's allure\", \"Electrifying graphics and vintage look\"),\n    # \"What we don't like\" bullets\n    @(\"High level of volatility may not appeal to everyone\", \"High level of volatility\"),\n    @(\"30 paylines may seem limited compared to some other online slots\", \"Could offer more paylines\"),\n    # Meta description (italic run)\n    @(\"Discover Krakatoa Lucky Break, an online slot with frequent wins and huge jackpots. Play for free and read our review of the game's mechanics.\", \"Read our review of Krakatoa Lucky Break, a slot game with high volatility and big jackpots. Play for free!\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n}\n"}
